# Updates cryptos list values (price + 1h volume change) to match the
# "Updated cryptos list on Mon May 29 15:50:14 UTC 2023 with GitHub Actions" commit.
# Column D ("Price") cells hold text-formatted numbers (e.g. "317.00", "27.893.59")
# so a leading apostrophe forces Excel to keep them as text instead of parsing them
# as numeric values (which would also strip formatting like trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.879.03"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "'1.905.91"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'316.93"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'0.4837"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("D8").Value = "'0.3797"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.07371"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'0.9311"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'20.77"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "'0.07748"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "'1.895.03"
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").Value = "'5.484"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "'6.633"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'91.84"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D18").Value = "'0.000008863"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D20").Value = "'27.943.81"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "'14.64"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'5.159"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "'2.161.53"
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "'1.921"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'154.88"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'18.47"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'2.136"
$ws.Range("E28").Value = "  +6.72%  "
$ws.Range("D29").Value = "'117.31"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").Value = "'4.959"
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "'0.08975"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'3.238"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("D33").Value = "'1.254"
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("D34").Value = "'0.7667"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "'4.660"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'0.02047"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'2.537"
$ws.Range("E37").Value = "  -6.14%  "
$ws.Range("D38").Value = "'1.096"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05288"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.998"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "'0.5477"
$ws.Range("D42").Value = "'6.960"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "'0.1527"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "'8.408"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'109.23"
$ws.Range("E45").Value = "  +5.15%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.62"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'0.4810"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'1.652"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'67.65"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "'0.06088"
$ws.Range("E51").Value = "  -0.17%  "
